# Apply edits described by the commit diff:
# - Add "UserType"/"Staff" column to ManageUsers sheet (C1/C2)
# - Add "Dropdownvalue"/"Toys" column to Subcategory sheet (B1/B2)
# - Change active sheet/tab selection from Category to Subcategory
# - Update selection/active cells on the affected sheets

$wb = $excel.ActiveWorkbook

$wsSubcategory  = $wb.Worksheets.Item("Subcategory")
$wsManageUsers  = $wb.Worksheets.Item("ManageUsers")

# --- ManageUsers sheet (sheet5): add new column C with header + value ---
$wsManageUsers.Range("C1").Value = "UserType"
$wsManageUsers.Range("C2").Value = "Staff"
$wsManageUsers.Range("C2").Select()

# --- Subcategory sheet (sheet4): add new column B with header + value ---
$wsSubcategory.Range("B1").Value = "Dropdownvalue"
$wsSubcategory.Range("B2").Value = "Toys"
$wsSubcategory.Range("A7").Select()

# --- Tab selection: Category is no longer the tab-selected sheet, Subcategory is ---
$wsSubcategory.Activate()
$wsSubcategory.Select()

$wb.Save()
